$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 538 (old rows 538-547 shift down to 540-549)
$ws.Rows.Item(538).Insert()
$ws.Rows.Item(538).Insert()

# New row 538
$ws.Range("A538").Value = 3
$ws.Range("B538").Value = "Femacal de La Calera"
$ws.Range("C538").Value = "Coquimbo"
$ws.Range("D538").Value = 44656
$ws.Range("E538").Value = 5
$ws.Range("F538").Value = 100114001
$ws.Range("G538").Value = "Papa"
$ws.Range("H538").Value = "Asterix"
$ws.Range("I538").Value = "1a (cosecha)"
$ws.Range("J538").Value = 180
$ws.Range("K538").Value = 7000
$ws.Range("L538").Value = 7000
$ws.Range("M538").Value = 7000
$ws.Range("N538").Value = '$/saco 25 kilos'
$ws.Range("O538").Value = "Provincia de Talca"
$ws.Range("P538").Value = 280
$ws.Range("Q538").Value = 25
$ws.Range("R538").Value = "Hortaliza"

# New row 539
$ws.Range("A539").Value = 3
$ws.Range("B539").Value = "Femacal de La Calera"
$ws.Range("C539").Value = "Coquimbo"
$ws.Range("D539").Value = 44656
$ws.Range("E539").Value = 5
$ws.Range("F539").Value = 100114001
$ws.Range("G539").Value = "Papa"
$ws.Range("H539").Value = "Rosara"
$ws.Range("I539").Value = "1a (cosecha)"
$ws.Range("J539").Value = 310
$ws.Range("K539").Value = 6500
$ws.Range("L539").Value = 7000
$ws.Range("M539").Value = 6742
$ws.Range("N539").Value = '$/saco 25 kilos'
$ws.Range("O539").Value = "Provincia de Talca"
$ws.Range("P539").Value = 270
$ws.Range("Q539").Value = 25
$ws.Range("R539").Value = "Hortaliza"
